# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) per item in each crafting-leve profit sheet, sourced from the
# latest Universalis market data snapshot.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 11: "Gotta Bounce" / Rubber
$ws.Range("H11").Value = 87.818184
$ws.Range("I11").Value = 87.818184
$ws.Range("K11").Value = 87.818184
$ws.Range("M11").Value = 52.181816
# Row 18: "You Grow, Girl" / Growth Formula Beta
$ws.Range("H18").Value = 1198
$ws.Range("I18").Value = 1198
$ws.Range("K18").Value = 1198
$ws.Range("M18").Value = -914
# Row 38: "Just Give Him a Serum" / Hi-Potion of Strength
$ws.Range("H38").Value = 437.81818
$ws.Range("I38").Value = 313.44446
$ws.Range("J38").Value = 997.5
$ws.Range("K38").Value = 940.33338
$ws.Range("L38").Value = 2992.5
$ws.Range("M38").Value = -568.33338
$ws.Range("N38").Value = -3736.5
# Row 40: "Stuck in the Moment" / Horn Glue
$ws.Range("H40").Value = 3357.0715
$ws.Range("I40").Value = 6224.75
$ws.Range("J40").Value = 2210
$ws.Range("K40").Value = 6224.75
$ws.Range("L40").Value = 2210
$ws.Range("M40").Value = -6049.75
$ws.Range("N40").Value = -2560
# Row 74: "Adhesive of Antipathy" / Wing Glue
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77: "It's Gonna Grow Back (L)" / Wing Glue
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 105: "Ultimate Official Strategy Guide" / Gazelleskin Codex
$ws.Range("H105").Value = 55277.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 55277.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 55277.5
$ws.Range("N105").Value = -62265.5
$ws.Range("M105").ClearContents()
# Row 117: "A Greater Grimoire" / Zonureskin Grimoire
$ws.Range("H117").Value = 94500
$ws.Range("J117").Value = 94500
$ws.Range("L117").Value = 94500
$ws.Range("N117").Value = -103678

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2: "Ain't Got No Ingots" / Bronze Ingot
$ws.Range("H2").Value = 1284.5834
$ws.Range("I2").Value = 579.44446
$ws.Range("J2").Value = 3400
$ws.Range("K2").Value = 579.44446
$ws.Range("L2").Value = 3400
$ws.Range("M2").Value = -466.44446
$ws.Range("N2").Value = -3626
# Row 74: "As the Bolt Flies" / Titanium Nugget
$ws.Range("H74").Value = 500.52942
$ws.Range("I74").Value = 500.52942
$ws.Range("K74").Value = 500.52942
$ws.Range("M74").Value = 373.47058
# Row 77: "Heavy Metal Banned (L)" / Titanium Nugget
$ws.Range("H77").Value = 500.52942
$ws.Range("I77").Value = 500.52942
$ws.Range("K77").Value = 2502.6471
$ws.Range("M77").Value = 1865.3529
# Row 80: "A Squire to Inspire" / Titanium Hoplon
$ws.Range("H80").Value = 99997
$ws.Range("J80").Value = 99997
$ws.Range("L80").Value = 99997
$ws.Range("N80").Value = -101993
# Row 83: "All's Fair in Highborn Assassination (L)" / Titanium Hoplon
$ws.Range("H83").Value = 99997
$ws.Range("J83").Value = 99997
$ws.Range("L83").Value = 299991
$ws.Range("N83").Value = -309975
# Row 95: "Shielded Life" / High Steel Scutum
$ws.Range("H95").Value = 30449.5
$ws.Range("J95").Value = 30449.5
$ws.Range("L95").Value = 30449.5
$ws.Range("N95").Value = -35941.5
# Row 97: "Ore for Me" / High Steel Ingot
$ws.Range("H97").Value = 577
$ws.Range("I97").Value = 555
$ws.Range("J97").Value = 786
$ws.Range("K97").Value = 555
$ws.Range("L97").Value = 786
$ws.Range("M97").Value = -59
$ws.Range("N97").Value = -1778
# Row 110: "Scheduled Maintenance" / Deepgold Ingot
$ws.Range("H110").Value = 1877.9546
$ws.Range("I110").Value = 1974.4736
$ws.Range("J110").Value = 1266.6666
$ws.Range("K110").Value = 1974.4736
$ws.Range("L110").Value = 1266.6666
$ws.Range("M110").Value = 70.52639999999997
$ws.Range("N110").Value = -5356.6666
# Row 116: "No Scope" / Titanbronze Ingot
$ws.Range("H116").Value = 1284.5834
$ws.Range("I116").Value = 579.44446
$ws.Range("J116").Value = 3400
$ws.Range("K116").Value = 579.44446
$ws.Range("L116").Value = 3400
$ws.Range("M116").Value = 1714.55554
$ws.Range("N116").Value = -7988

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3: "Hells Bells" / Bronze Ingot
$ws.Range("H3").Value = 1284.5834
$ws.Range("I3").Value = 579.44446
$ws.Range("J3").Value = 3400
$ws.Range("K3").Value = 579.44446
$ws.Range("L3").Value = 3400
$ws.Range("M3").Value = -465.44446
$ws.Range("N3").Value = -3628
# Row 107: "The Gold Experience" / Deepgold Nugget
$ws.Range("H107").Value = 2107.2
$ws.Range("I107").Value = 2063.5557
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 2063.5557
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -143.5556999999999
$ws.Range("N107").Value = -6340

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16: "Raise the Roof" / Ash Lumber
$ws.Range("H16").Value = 1088.625
$ws.Range("I16").Value = 672.8571
$ws.Range("K16").Value = 672.8571
$ws.Range("M16").Value = -385.8570999999999
# Row 31: "Wall Not Found" / Walnut Lumber
$ws.Range("H31").Value = 2578.9722
$ws.Range("I31").Value = 1946
$ws.Range("K31").Value = 1946
$ws.Range("M31").Value = -1651
# Row 34: "Armoires of the Rich and Famous" / Walnut Lumber
$ws.Range("H34").Value = 2578.9722
$ws.Range("I34").Value = 1946
$ws.Range("K34").Value = 1946
$ws.Range("M34").Value = -1744
# Row 41: "The Lone Bowman" / Oak Longbow
$ws.Range("H41").Value = 16756.584
$ws.Range("I41").Value = 4093.1667
$ws.Range("K41").Value = 4093.1667
$ws.Range("M41").Value = -3665.1667
# Row 50: "The Arsenal of Theocracy" / Cobalt Halberd
$ws.Range("H50").Value = 22833
$ws.Range("J50").Value = 22833
$ws.Range("L50").Value = 22833
$ws.Range("N50").Value = -24083
# Row 113: "Patient Patients" / White Ash Lumber
$ws.Range("H113").Value = 1088.625
$ws.Range("I113").Value = 672.8571
$ws.Range("K113").Value = 672.8571
$ws.Range("M113").Value = 1497.1429
# Row 134: "Wood You Be Quiet" / Ceiba Lumber
$ws.Range("H134").Value = 2511.5715
$ws.Range("I134").Value = 2334.6843
$ws.Range("K134").Value = 7004.0529
$ws.Range("M134").Value = -4469.0529

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 29: "For Crumbs' Sake" / Honey Muffin
$ws.Range("H29").Value = 333525.34
$ws.Range("I29").Value = 1000014.5
$ws.Range("J29").Value = 280.75
$ws.Range("K29").Value = 3000043.5
$ws.Range("L29").Value = 842.25
$ws.Range("M29").Value = -2999766.5
$ws.Range("N29").Value = -1396.25
# Row 63: "The Next to Last Supper" / Stuffed Cabbage Rolls
$ws.Range("H63").Value = 3431.6667
$ws.Range("I63").Value = 3447.5
$ws.Range("K63").Value = 10342.5
$ws.Range("M63").Value = -9593.5
# Row 66: "Nostalgia through the Stomach (L)" / Stuffed Cabbage Rolls
$ws.Range("H66").Value = 3431.6667
$ws.Range("I66").Value = 3447.5
$ws.Range("K66").Value = 31027.5
$ws.Range("M66").Value = -27283.5
# Row 68: "Such a Butter Face" / Fermented Butter
$ws.Range("H68").Value = 3026.5557
$ws.Range("I68").Value = 2741.1667
$ws.Range("J68").Value = 3597.3333
$ws.Range("K68").Value = 8223.500100000001
$ws.Range("L68").Value = 10791.9999
$ws.Range("M68").Value = -7412.500100000001
$ws.Range("N68").Value = -12413.9999
# Row 69: "Loving That Muffin Top" / Ishgardian Muffin
$ws.Range("H69").Value = 2671.2856
$ws.Range("J69").Value = 2979.8
$ws.Range("L69").Value = 8939.400000000001
$ws.Range("N69").Value = -10561.4
# Row 71: "No Margarine of Error (L)" / Fermented Butter
$ws.Range("H71").Value = 3026.5557
$ws.Range("I71").Value = 2741.1667
$ws.Range("J71").Value = 3597.3333
$ws.Range("K71").Value = 24670.5003
$ws.Range("L71").Value = 32375.9997
$ws.Range("M71").Value = -20614.5003
$ws.Range("N71").Value = -40487.9997
# Row 72: "Muffin of the Morn (L)" / Ishgardian Muffin
$ws.Range("H72").Value = 2671.2856
$ws.Range("J72").Value = 2979.8
$ws.Range("L72").Value = 26818.2
$ws.Range("N72").Value = -34930.2
# Row 80: "Saucy for a Suitor" / Hollandaise Sauce
$ws.Range("H80").Value = 2733
$ws.Range("J80").Value = 2733
$ws.Range("L80").Value = 8199
$ws.Range("N80").Value = -10071
# Row 81: "It Goes Down Smoothly" / Frozen Spirits
$ws.Range("H81").Value = 3100
$ws.Range("I81").Value = 1500
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 4500
$ws.Range("L81").Value = 10500
$ws.Range("M81").Value = -3377
$ws.Range("N81").Value = -12746
# Row 83: "Saved by the Sauce (L)" / Hollandaise Sauce
$ws.Range("H83").Value = 2733
$ws.Range("J83").Value = 2733
$ws.Range("L83").Value = 24597
$ws.Range("N83").Value = -33957
# Row 84: "Quenching the Flame (L)" / Frozen Spirits
$ws.Range("H84").Value = 3100
$ws.Range("I84").Value = 1500
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 13500
$ws.Range("L84").Value = 31500
$ws.Range("M84").Value = -7884
$ws.Range("N84").Value = -42732
# Row 137: "Creative Chocolate" / Gateau au Chocolat
$ws.Range("H137").Value = 4503.3335
$ws.Range("I137").Value = 5132.5
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 15397.5
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -10297.5
$ws.Range("N137").Value = -22200
# Row 138: "Bring Me Your Tacos" / Tacos Al Pastor
$ws.Range("H138").Value = 5894.5
$ws.Range("I138").Value = 5894.5
$ws.Range("K138").Value = 17683.5
$ws.Range("M138").Value = -12543.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97: "If I'd a Koppranickel for Every Time..." / Koppranickel Ingot
$ws.Range("H97").Value = 815.25
$ws.Range("I97").Value = 815.25
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 815.25
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -319.25
$ws.Range("N97").ClearContents()
# Row 98: "Cutting Deals" / Durium Smallsword
$ws.Range("H98").Value = 12216.333
$ws.Range("J98").Value = 12216.333
$ws.Range("L98").Value = 12216.333
$ws.Range("N98").Value = -18206.333
# Row 101: "Best-laid Planispheres" / Dual-plated Durium Planisphere
$ws.Range("H101").Value = 39425.5
$ws.Range("J101").Value = 39425.5
$ws.Range("L101").Value = 39425.5
$ws.Range("N101").Value = -45915.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7: "Tan Before the Ban" / Leather
$ws.Range("H7").Value = 1623.375
$ws.Range("I7").Value = 1398.8
$ws.Range("J7").Value = 1997.6666
$ws.Range("K7").Value = 1398.8
$ws.Range("L7").Value = 1997.6666
$ws.Range("M7").Value = -1286.8
$ws.Range("N7").Value = -2221.6666
# Row 46: "Supply Side Logic" / Boar Leather
$ws.Range("H46").Value = 500400
$ws.Range("I46").Value = 801
$ws.Range("J46").Value = 999999
$ws.Range("K46").Value = 801
$ws.Range("L46").Value = 999999
$ws.Range("M46").Value = -613
$ws.Range("N46").Value = -1000375
# Row 126: "Battered Books" / Saiga Leather
$ws.Range("H126").Value = 1623.375
$ws.Range("I126").Value = 1398.8
$ws.Range("J126").Value = 1997.6666
$ws.Range("K126").Value = 4196.4
$ws.Range("L126").Value = 5992.9998
$ws.Range("M126").Value = -1726.4
$ws.Range("N126").Value = -10932.9998
# Row 136: "Respect for Br'aax" / Br'aax Leather
$ws.Range("H136").Value = 6627.6665
$ws.Range("I136").Value = 5497.5
$ws.Range("K136").Value = 16492.5
$ws.Range("M136").Value = -13942.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 117: "The Hunt Continues" / Ovim Wool Muffed Met of Casting
$ws.Range("H117").Value = 97704.5
$ws.Range("J117").Value = 97704.5
$ws.Range("L117").Value = 97704.5
$ws.Range("N117").Value = -106882.5
